# БюджетДекабрь.xlsx — "Update Exec13 upr 1, 2, 4"
# Adds a salary (ЗП) income row and two new expense rows (Ипотека / Доскрочка),
# moves the "Подарки НГ" label from N4 into the newly merged M4:N4 header cell,
# widens column C and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: new income entry (ЗП) in the "Поступления" block -------------
$ws.Range("E5").Copy()
$ws.Range("A6").PasteSpecial(-4122)     # xlPasteFormats - reuse the date style (s=2)
$ws.Range("A6").Value = 44924           # 29-Dec-2022
$ws.Range("B6").Value = "ЗП"
$ws.Range("C6").Value = 41371.86

# --- Row 7: new expense entry (Ипотека) in the "Кредиты" block -----------
$ws.Range("I5").Copy()
$ws.Range("I7").PasteSpecial(-4122)     # xlPasteFormats - reuse the date style (s=2)
$ws.Range("I7").Value = 44924           # 29-Dec-2022
$ws.Range("J7").Value = "Ипотека"
$ws.Range("K7").Value = 16561

# --- Row 8: new expense entry (Доскрочка) in the "Кредиты" block ---------
$ws.Range("J8").Value = "Доскрочка"
$ws.Range("K8").Value = 15000

# --- Move the "Подарки НГ" header label from N4 into M4, merge M4:N4 -----
$ws.Range("I4:J4").Copy()
$ws.Range("M4:N4").PasteSpecial(-4122)  # xlPasteFormats - reuse the merged-header style
$label = $ws.Range("N4").Value()
$ws.Range("M4").Value = $label
$ws.Range("N4").ClearContents()
$ws.Range("M4:N4").Merge()

# --- Column C gets wider (6.43 -> ~10.71 chars) ---------------------------
$ws.Columns("C").ColumnWidth = 9.9

# --- Active selection moves to N11 ----------------------------------------
$ws.Range("N11").Select()
